# "Update LichLamViec 26/1 - 7/2"
# Adds two new weekly blocks (Tuan: 16, Tuan: 17) to the "Lich lam viec" sheet
# and extends the note for the last day of "Tuan: 15".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Tuan 15 (rows 108-111): extend Friday's note and fill in the previously
#    empty Saturday cell.
# ---------------------------------------------------------------------------
$ws.Range("D119").Value = "Tiềm hiểu về cách cào dữ liệu từ trang web khác bằng cách đọc html. Tìm hiểu về CrystalReport"
$ws.Range("E119").Value = "Nghỉ"
$ws.Rows(119).RowHeight = 63.75

# ---------------------------------------------------------------------------
# Helper block: rows 108-111 hold one fully-formatted "week" template
# (header / dates-formula / weekday-formula / notes). We reuse its
# formatting for the two new blocks below so fonts/fills/borders match
# exactly, then overwrite values & formulas per destination week.
# ---------------------------------------------------------------------------

function Copy-WeekFormat($destHeaderRow) {
    $h = $destHeaderRow        # header row (e.g. 123)      <- template row 108
    $d = $destHeaderRow + 1    # dates-formula row            <- template row 109
    $w = $destHeaderRow + 2    # weekday-formula row           <- template row 110
    $n = $destHeaderRow + 3    # notes row                      <- template row 111

    $ws.Range("A108:S108").Copy() | Out-Null
    $ws.Range("A$h").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Range("A109:I109").Copy() | Out-Null
    $ws.Range("A$d").PasteSpecial(-4122) | Out-Null

    $ws.Range("A110:I110").Copy() | Out-Null
    $ws.Range("A$w").PasteSpecial(-4122) | Out-Null

    $ws.Range("A111:I111").Copy() | Out-Null
    $ws.Range("A$n").PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) Tuan 16 block: rows 123 (header), 124 (dates), 125 (weekday names),
#    126 (notes - everyone off, 12.75pt tall).
# ---------------------------------------------------------------------------
Copy-WeekFormat 123

$ws.Range("B123").Value = "Tuần: 16`n"
$ws.Range("C123").Value = 45683

$ws.Range("C124").Formula = "=C123"
$ws.Range("D124").Formula = "=C123+1"
$ws.Range("E124").Formula = "=C123+2"
$ws.Range("F124").Formula = "=C123+3"
$ws.Range("G124").Formula = "=C123+4"
$ws.Range("H124").Formula = "=C123+5"
$ws.Range("I124").Formula = "=C123+6"

$ws.Range("C125").Formula = '=UPPER(TEXT(C124, "DDDD"))'
$ws.Range("D125").Formula = '=UPPER(TEXT(D124, "DDDD"))'
$ws.Range("E125").Formula = '=UPPER(TEXT(E124, "DDDD"))'
$ws.Range("F125").Formula = '=UPPER(TEXT(F124, "DDDD"))'
$ws.Range("G125").Formula = '=UPPER(TEXT(G124, "DDDD"))'
$ws.Range("H125").Formula = '=UPPER(TEXT(H124, "DDDD"))'
$ws.Range("I125").Formula = '=UPPER(TEXT(I124, "DDDD"))'

$ws.Range("C126").Value = "Nghỉ"
$ws.Range("D126").Value = "Nghỉ"
$ws.Range("E126").Value = "Nghỉ"
$ws.Range("F126").Value = "Nghỉ"
$ws.Range("G126").Value = "Nghỉ"
$ws.Range("H126").Value = "Nghỉ"
$ws.Range("I126").Value = "Nghỉ"

$ws.Rows(123).RowHeight = 30
$ws.Rows(124).RowHeight = 15.75
$ws.Rows(125).RowHeight = 15.75
$ws.Rows(126).RowHeight = 12.75

$ws.Range("E123:I123").Merge()
$ws.Range("J123:N123").Merge()
$ws.Range("O123:S123").Merge()

# ---------------------------------------------------------------------------
# 3) Tuan 17 block: rows 130 (header), 131 (dates), 132 (weekday names),
#    133 (notes about Vue.js).
# ---------------------------------------------------------------------------
Copy-WeekFormat 130

$ws.Range("B130").Value = "Tuần: 17`n"
$ws.Range("C130").Value = 45690

$ws.Range("C131").Formula = "=C130"
$ws.Range("D131").Formula = "=C130+1"
$ws.Range("E131").Formula = "=C130+2"
$ws.Range("F131").Formula = "=C130+3"
$ws.Range("G131").Formula = "=C130+4"
$ws.Range("H131").Formula = "=C130+5"
$ws.Range("I131").Formula = "=C130+6"

$ws.Range("C132").Formula = '=UPPER(TEXT(C131, "DDDD"))'
$ws.Range("D132").Formula = '=UPPER(TEXT(D131, "DDDD"))'
$ws.Range("E132").Formula = '=UPPER(TEXT(E131, "DDDD"))'
$ws.Range("F132").Formula = '=UPPER(TEXT(F131, "DDDD"))'
$ws.Range("G132").Formula = '=UPPER(TEXT(G131, "DDDD"))'
$ws.Range("H132").Formula = '=UPPER(TEXT(H131, "DDDD"))'
$ws.Range("I132").Formula = '=UPPER(TEXT(I131, "DDDD"))'

$ws.Range("C133").Value = "Nghỉ"
$ws.Range("D133").Value = "Tìm hiểu về Vue.js."
$ws.Range("E133").Value = "Tìm hiểu về Vue.js. Về component "
$ws.Range("F133").Value = "Tìm hiểu về Vue.js. Về component"
$ws.Range("G133").Value = "Tìm hiểu về Vue.js, Về hiệu ứng chuyển động"
$ws.Range("H133").Value = "Tìm hiểu về Vue.js. Chỉnh lại form SendMessToAsset chọn địa chỉ BootstrapServer (FOX, LLQ)"
$ws.Range("I133").Value = "Nghỉ"

$ws.Rows(130).RowHeight = 30
$ws.Rows(131).RowHeight = 15.75
$ws.Rows(132).RowHeight = 15.75
$ws.Rows(133).RowHeight = 63.75

$ws.Range("E130:I130").Merge()
$ws.Range("J130:N130").Merge()
$ws.Range("O130:S130").Merge()

# ---------------------------------------------------------------------------
# 4) Restore blank-row heights below each new block that PasteSpecial may
#    not have touched, and keep the plain filler rows beneath untouched.
# ---------------------------------------------------------------------------
for ($r = 127; $r -le 129; $r++) { $ws.Rows($r).RowHeight = 15.75 }
for ($r = 134; $r -le 144; $r++) { $ws.Rows($r).RowHeight = 15.75 }

# ---------------------------------------------------------------------------
# 5) View state (scroll position / active cell) - cosmetic but mirrors diff.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 115
$ws.Range("G135").Select()
